# Welcome-mail generator: turn Rishabh's plain-text username into a
# real e-mail address and wire it up as a clickable mailto: hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix up the e-mail column value for Rishabh (row 2) so it is a
#    real address instead of a bare username.
$emailCell = $ws.Range("B2")
$emailCell.Value = "ridhabhthaney@gmail.com"

# 2. Turn that value into a mailto: hyperlink so the generated welcome
#    mail can be sent with one click.
$ws.Hyperlinks.Add($emailCell, "mailto:ridhabhthaney@gmail.com", "", "", "ridhabhthaney@gmail.com")

# 3. Nudge row 2 a touch taller and widen the columns slightly so the
#    new hyperlink text is comfortably readable.
$ws.Rows.Item(2).RowHeight = 15.7
$ws.Columns.Item(1).ColumnWidth = 10.1666666666667
$ws.Columns.Item(2).ColumnWidth = 33.5
$ws.Columns.Item(3).ColumnWidth = 17.5
$ws.Columns.Item(4).ColumnWidth = 10.1666666666667

# 4. Leave the selection where the author ended up while working on
#    this (row 11 of the sheet).
$null = $ws.Range("B11").Select()
